$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "InitialScreening"
$ws.Range("B3").Value = "SD\CA.xlsx"
$ws.Range("C3").Value = "InitialScreening"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B4").Select()
